# Applies the updated loading_percent results for "case with 380 kV done"
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("B2").Value = 7.773726052723036
$ws.Range("D2").Value = 3.229918762227883
$ws.Range("E2").Value = 10.85449565201054
$ws.Range("F2").Value = 27.32259848452737
$ws.Range("G2").Value = 38.37620787344353
$ws.Range("H2").Value = 13.91355051060317
$ws.Range("M2").Value = 19.29931131823499
$ws.Range("N2").Value = 17.30400529927183
$ws.Range("B3").Value = 7.699769585378305
$ws.Range("D3").Value = 3.217610194302631
$ws.Range("E3").Value = 10.95388439539627
$ws.Range("F3").Value = 26.36821196932372
$ws.Range("G3").Value = 36.44911460145707
$ws.Range("H3").Value = 13.70211742965484
$ws.Range("M3").Value = 18.45188709384038
$ws.Range("N3").Value = 17.26266740747192
$ws.Range("B4").Value = 7.65580617498127
$ws.Range("D4").Value = 3.211470150608485
$ws.Range("E4").Value = 11.01793081773803
$ws.Range("F4").Value = 25.77620538914855
$ws.Range("G4").Value = 35.22847362765715
$ws.Range("H4").Value = 13.57606642251865
$ws.Range("M4").Value = 17.91451140887566
$ws.Range("N4").Value = 17.23980560046354
$ws.Range("B5").Value = 7.638274756851349
$ws.Range("D5").Value = 3.209323052025934
$ws.Range("E5").Value = 11.0447934376487
$ws.Range("F5").Value = 25.53388664156654
$ws.Range("G5").Value = 34.72241638740584
$ws.Range("H5").Value = 13.52571201964763
$ws.Range("M5").Value = 17.69155044316827
$ws.Range("N5").Value = 17.23113123132711
$ws.Range("B6").Value = 7.635387454430822
$ws.Range("D6").Value = 3.208987916862917
$ws.Range("E6").Value = 11.04930015029949
$ws.Range("F6").Value = 25.49359856675179
$ws.Range("G6").Value = 34.63788911269963
$ws.Range("H6").Value = 13.51741368993776
$ws.Range("M6").Value = 17.65429774929389
$ws.Range("N6").Value = 17.22972987062609
$ws.Range("B7").Value = 7.655568158807616
$ws.Range("D7").Value = 3.2114397589402
$ws.Range("E7").Value = 11.0182900017394
$ws.Range("F7").Value = 25.77294112675085
$ws.Range("G7").Value = 35.22168262945369
$ws.Range("H7").Value = 13.57538314289476
$ws.Range("M7").Value = 17.91152013479003
$ws.Range("N7").Value = 17.23968600467914
$ws.Range("B8").Value = 7.747938144797139
$ws.Range("D8").Value = 3.225379908084616
$ws.Range("E8").Value = 10.88813989565872
$ws.Range("F8").Value = 26.99502717524225
$ws.Range("G8").Value = 37.71993346052442
$ws.Range("H8").Value = 13.83990537445398
$ws.Range("M8").Value = 19.01083369296762
$ws.Range("N8").Value = 17.28923414937545
$ws.Range("B9").Value = 7.939624014049664
$ws.Range("D9").Value = 3.26399729697104
$ws.Range("E9").Value = 10.65673097279664
$ws.Range("F9").Value = 29.32621113141824
$ws.Range("G9").Value = 42.29410216595762
$ws.Range("H9").Value = 14.38560804881461
$ws.Range("M9").Value = 21.01956363572681
$ws.Range("N9").Value = 17.40600716059572
$ws.Range("B10").Value = 8.085612679638324
$ws.Range("D10").Value = 3.299251211636756
$ws.Range("E10").Value = 10.50101102537828
$ws.Range("F10").Value = 30.97794355530005
$ws.Range("G10").Value = 45.42582351448262
$ws.Range("H10").Value = 14.7989628869554
$ws.Range("M10").Value = 22.3926124619051
$ws.Range("N10").Value = 17.50321596352093
$ws.Range("B11").Value = 8.152870062141595
$ws.Range("D11").Value = 3.331426226037686
$ws.Range("E11").Value = 10.43322635527094
$ws.Range("F11").Value = 31.7123434070001
$ws.Range("G11").Value = 46.79597674902132
$ws.Range("H11").Value = 14.9888489031292
$ws.Range("M11").Value = 22.99279631955855
$ws.Range("N11").Value = 17.54979497354602
$ws.Range("B12").Value = 8.178437153016974
$ws.Range("D12").Value = 3.355453717100469
$ws.Range("E12").Value = 10.40799335911325
$ws.Range("F12").Value = 31.98773081461442
$ws.Range("G12").Value = 47.30667104756485
$ws.Range("H12").Value = 15.06094899093316
$ws.Range("M12").Value = 23.21641600497285
$ws.Range("N12").Value = 17.56776139477675
$ws.Range("B13").Value = 8.172926829097879
$ws.Range("D13").Value = 3.350295205885344
$ws.Range("E13").Value = 10.4134084184197
$ws.Range("F13").Value = 31.92854590226108
$ws.Range("G13").Value = 47.19705096252505
$ws.Range("H13").Value = 15.04541342281846
$ws.Range("M13").Value = 23.16842020911741
$ws.Range("N13").Value = 17.56387759931358
$ws.Range("B14").Value = 8.154971648708047
$ws.Range("D14").Value = 3.33341029080885
$ws.Range("E14").Value = 10.43114171146295
$ws.Range("F14").Value = 31.73505552958859
$ws.Range("G14").Value = 46.83815689434114
$ws.Range("H14").Value = 14.99477712889079
$ws.Range("M14").Value = 23.01126768341198
$ws.Range("N14").Value = 17.55126657402592
$ws.Range("B15").Value = 8.143985672118381
$ws.Range("D15").Value = 3.32302034836781
$ws.Range("E15").Value = 10.44206048930841
$ws.Range("F15").Value = 31.61617613536437
$ws.Range("G15").Value = 46.61725368891234
$ws.Range("H15").Value = 14.96378415577902
$ws.Range("M15").Value = 22.91452708182083
$ws.Range("N15").Value = 17.54358433063172
$ws.Range("B16").Value = 8.081232269048323
$ws.Range("D16").Value = 3.298135900937261
$ws.Range("E16").Value = 10.50550205504354
$ws.Range("F16").Value = 30.92958310534992
$ws.Range("G16").Value = 45.3351558827819
$ws.Range("H16").Value = 14.78658444683272
$ws.Range("M16").Value = 22.35288465420948
$ws.Range("N16").Value = 17.50021842375544
$ws.Range("B17").Value = 8.042935097681069
$ws.Range("D17").Value = 3.288527236950203
$ws.Range("E17").Value = 10.54520092735451
$ws.Range("F17").Value = 30.50382738959005
$ws.Range("G17").Value = 44.53443055579642
$ws.Range("H17").Value = 14.67830036875361
$ws.Range("M17").Value = 22.00196918617873
$ws.Range("N17").Value = 17.47421101742254
$ws.Range("B18").Value = 8.020988650515662
$ws.Range("D18").Value = 3.28314034847426
$ws.Range("E18").Value = 10.56832224915518
$ws.Range("F18").Value = 30.25735973383091
$ws.Range("G18").Value = 44.06876135025603
$ws.Range("H18").Value = 14.61619645190364
$ws.Range("M18").Value = 21.79784013636036
$ws.Range("N18").Value = 17.45947490265145
$ws.Range("B19").Value = 8.013572603771969
$ws.Range("D19").Value = 3.281340496897868
$ws.Range("E19").Value = 10.57620022843673
$ws.Range("F19").Value = 30.17364674215251
$ws.Range("G19").Value = 43.91022624117641
$ws.Range("H19").Value = 14.59520194480194
$ws.Range("M19").Value = 21.72833674274569
$ws.Range("N19").Value = 17.45452408334116
$ws.Range("B20").Value = 8.047003677077717
$ws.Range("D20").Value = 3.289535643757433
$ws.Range("E20").Value = 10.54094517608986
$ws.Range("F20").Value = 30.54931584578891
$ws.Range("G20").Value = 44.6202008204185
$ws.Range("H20").Value = 14.68980947130236
$ws.Range("M20").Value = 22.03956304908913
$ws.Range("N20").Value = 17.47695658576031
$ws.Range("B21").Value = 8.160243042688977
$ws.Range("D21").Value = 3.338379698192915
$ws.Range("E21").Value = 10.42592121923824
$ws.Range("F21").Value = 31.7919640058358
$ws.Range("G21").Value = 46.94379630497381
$ws.Range("H21").Value = 15.00964552222715
$ws.Range("M21").Value = 23.05752747733077
$ws.Range("N21").Value = 17.55496192754454
$ws.Range("B22").Value = 8.234812993339558
$ws.Range("D22").Value = 3.407632754068791
$ws.Range("E22").Value = 10.35328374262822
$ws.Range("F22").Value = 32.58818374093756
$ws.Range("G22").Value = 48.4147695470712
$ws.Range("H22").Value = 15.21977793046083
$ws.Range("M22").Value = 23.70145932136301
$ws.Range("N22").Value = 17.60784914112065
$ws.Range("B23").Value = 8.194969963266729
$ws.Range("D23").Value = 3.370866913201756
$ws.Range("E23").Value = 10.39182070684568
$ws.Range("F23").Value = 32.16476492662208
$ws.Range("G23").Value = 47.63413237992894
$ws.Range("H23").Value = 15.10754873506805
$ws.Range("M23").Value = 23.35977785218243
$ws.Range("N23").Value = 17.57945158709864
$ws.Range("B24").Value = 8.045164049064688
$ws.Range("D24").Value = 3.289079315357258
$ws.Range("E24").Value = 10.54286827304898
$ws.Range("F24").Value = 30.52875577338337
$ws.Range("G24").Value = 44.58144063446223
$ws.Range("H24").Value = 14.68460573325155
$ws.Range("M24").Value = 22.02257427516943
$ws.Range("N24").Value = 17.47571464102569
$ws.Range("B25").Value = 7.886768789509262
$ws.Range("D25").Value = 3.252339649641252
$ws.Range("E25").Value = 10.7168065810447
$ws.Range("F25").Value = 28.70487899956651
$ws.Range("G25").Value = 41.09497914901617
$ws.Range("H25").Value = 14.23550574225174
$ws.Range("M25").Value = 20.49337423980174
$ws.Range("N25").Value = 17.37237168578267
